# Insert a new weekly data row at row 39 (pushing existing rows 39-43 down
# to 40-44), then populate the new row 39 with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 39:43 down to 40:44, leaving a blank row 39 for the new entry.
$ws.Rows.Item(39).Insert()

# Fill in the new row 39 (carrying the same market/product metadata as the
# surrounding rows, with the new date and pricing figures).
$ws.Range("A39").Value = 5
$ws.Range("B39").Value = "Macroferia Regional de Talca"
$ws.Range("C39").Value = "Maule"
$ws.Range("D39").Value = 44491
$ws.Range("E39").Value = 7
$ws.Range("F39").Value = "Fruta"
$ws.Range("G39").Value = 100107
$ws.Range("H39").Value = "Otros"
$ws.Range("I39").Value = 100107002
$ws.Range("J39").Value = "Chirimoya"
$ws.Range("K39").Value = "Cultivar IV Región"
$ws.Range("L39").Value = "Primera"
$ws.Range("M39").Value = 100
$ws.Range("N39").Value = 25000
$ws.Range("O39").Value = 25000
$ws.Range("P39").Value = 25000
$ws.Range("Q39").Value = "$/bandeja 10 kilos"
$ws.Range("R39").Value = "Provincia de Limarí"
$ws.Range("S39").Value = 2500
$ws.Range("T39").Value = 10
